# Resultados INP preprocesadas con class weight.xlsx
# Script generico VGG19 / Xception / ResNet50V2 / ResNet101 / ResNet152 (Issue #35)
# Nuevas pruebas de otras redes de transfer learning; actualiza las metricas de la tabla
# y aplica el resaltado de color a la cabecera de la tabla.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the newly measured metrics (previously-empty cells) ---------------

# Baseline (row 4 -> VGG16): Acc / f1-score / AUC for Samsung
$ws.Range("B4").Value = 0.59
$ws.Range("C4").Value = 0.43
$ws.Range("D4").Value = 0.53

# Row 5 -> VGG19 (brand new row of results)
$ws.Range("B5").Value = 0.68
$ws.Range("C5").Value = 0.66
$ws.Range("D5").Value = 0.73
$ws.Range("E5").Value = 0.72
$ws.Range("F5").Value = 0.71
$ws.Range("G5").Value = 0.72

# Row 6 -> Exception (brand new row of results)
$ws.Range("B6").Value = 0.53
$ws.Range("C6").Value = 0.53
$ws.Range("D6").Value = 0.45
$ws.Range("E6").Value = 0.75
$ws.Range("F6").Value = 0.65
$ws.Range("G6").Value = 0.53

# Row 7 -> ResNet50V2
$ws.Range("B7").Value = 0.59
$ws.Range("C7").Value = 0.44
$ws.Range("D7").Value = 0.56000000000000005
$ws.Range("E7").Value = 0.74
$ws.Range("F7").Value = 0.65
$ws.Range("G7").Value = 0.53

# Row 8 -> ResNet101
$ws.Range("B8").Value = 0.6
$ws.Range("C8").Value = 0.46
$ws.Range("D8").Value = 0.64
$ws.Range("E8").Value = 0.75
$ws.Range("F8").Value = 0.65
$ws.Range("G8").Value = 0.64

# Row 9 -> ResNet152
$ws.Range("B9").Value = 0.6
$ws.Range("C9").Value = 0.46
$ws.Range("D9").Value = 0.66
$ws.Range("E9").Value = 0.75
$ws.Range("F9").Value = 0.65
$ws.Range("G9").Value = 0.65

# Row 10 (InceptionV3) stays blank for now.

# --- 2. Highlight the table header with colour fills -------------------------------
# Colours applied in the same order Excel records them in the MRU colour list
# (and therefore the same order they are appended to the fills table):
#   33CCFF, CCFFFF, FF5050, FF9999, FFFFCC

$ws.Range("B1:D1").Interior.Color = 16763955   # FF33CCFF - Samsung header
$ws.Range("B2:D2").Interior.Color = 16777164   # FFCCFFFF - Samsung sub-header
$ws.Range("E1:G1").Interior.Color = 5263615    # FFFF5050 - iPhone header
$ws.Range("E2:G2").Interior.Color = 10066431   # FFFF9999 - iPhone sub-header
$ws.Range("A3:G3").Interior.Color = 13434879   # FFFFFFCC - Baseline row

# --- 3. Move the active selection to G10, matching the saved view -----------------
$ws.Range("G10").Select()
